# ============================================================================
# feat: add 2022-Q4 data
#
# - Inserts a new "2022-Q4" worksheet right after "总计", before "2022-Q3"
#   (all later quarter sheets shift right by one position).
# - Populates "2022-Q4" with the fund-holding table for that quarter.
# - Updates the "总计" (summary) sheet with a new leading row for 2022-Q4
#   and renumbers/shifts the existing rows down.
# - Restores "2021-Q3" as the active tab (it was the active tab originally).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Borrow the header-row / index-column formatting (bold, bordered, centered
# style) from the existing "2022-Q3" sheet, which uses the identical layout.
$srcSheet = $wb.Worksheets.Item("2022-Q3")
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2").Copy()
$newSheet.Range("A2:A15").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund-holding rows for 2022-Q4.
# Columns: idx, code, name, size, stockPosition, positionPct, marketValue, rank
$fundData = @(
    @("0", "010591", "富国中国中小盘混合（QDII）美元", "32.21", "87.21", "2.13", "0.6861", "10"),
    @("1", "100061", "富国中国中小盘混合（QDII）人民币", "32.21", "87.21", "2.13", "0.6861", "10"),
    @("2", "005583", "易方达港股通红利灵活配置混合", "7.47", "88.74", "4.85", "0.3623", "3"),
    @("3", "012227", "景顺长城港股通全球竞争力混合A", "10.77", "85.21", "3.01", "0.3242", "9"),
    @("4", "011635", "富国港股通策略精选混合A", "7.15", "85.60", "1.29", "0.0922", "10"),
    @("5", "012228", "景顺长城港股通全球竞争力混合C", "1.65", "85.21", "3.01", "0.0497", "9"),
    @("6", "004266", "招商沪港深科技创新主题精选灵活配置混合A", "0.93", "90.63", "5.23", "0.0486", "2"),
    @("7", "005701", "上投摩根香港精选港股通混合A", "0.48", "89.99", "3.02", "0.0145", "10"),
    @("8", "010754", "招商沪港深科技创新主题精选灵活配置混合C", "0.25", "90.63", "5.23", "0.0131", "2"),
    @("9", "014146", "景顺长城港股通数字经济主题混合A", "0.60", "89.97", "2.04", "0.0122", "10"),
    @("10", "011636", "富国港股通策略精选混合C", "0.73", "85.60", "1.29", "0.0094", "10"),
    @("11", "005269", "华泰柏瑞港股通量化灵活配置混合", "0.54", "80.96", "1.61", "0.0087", "10"),
    @("12", "014147", "景顺长城港股通数字经济主题混合C", "0.29", "89.97", "2.04", "0.0059", "10"),
    @("13", "016921", "上投摩根香港精选港股通混合C", "0.02", "89.99", "3.02", "0.0006", "10")
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $rec = $fundData[$i]

    $newSheet.Cells.Item($r, 1).Value = [int]$rec[0]

    # Fund code, size, stock position, position pct and market value all look
    # numeric but are stored as TEXT in the source data, so force text via
    # NumberFormat before assignment, then strip the format back off (via
    # Style = "Normal") so the cell keeps the plain/default style seen
    # elsewhere in the sheet instead of picking up a "@"-formatted style.
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $rec[1]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $rec[2]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $rec[3]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $rec[4]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $rec[5]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $rec[6]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = [int]$rec[7]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add a new leading row for 2022-Q4
#    and shift the existing three rows down by one.
# ---------------------------------------------------------------------------

# A5 is a brand new cell; copy the bold/bordered index-column style from A4
# before writing into it so it matches the rest of column A.
$totalSheet.Cells.Item(4, 1).Copy()
$totalSheet.Cells.Item(5, 1).PasteSpecial(-4122)

# Write bottom-up so each source row is read before it gets overwritten.
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(5, 3).Value = 5
$totalSheet.Cells.Item(5, 4).Value = 1.11

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(4, 3).Value = 4
$totalSheet.Cells.Item(4, 4).Value = 2.18

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(3, 3).Value = 13
$totalSheet.Cells.Item(3, 4).Value = 2.46

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 14
$totalSheet.Cells.Item(2, 4).Value = 2.31

# ---------------------------------------------------------------------------
# 3. Restore "2021-Q3" as the selected/active tab (matches the original
#    workbook's selection state).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
